# Edit: add SPECIE_HUE gen row (replacing the empty EMPTY5 placeholder) and
# set its values to -1 for initial species, using the blend gen approach.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("initial_species")

# Row 12 is currently the "EMPTY5" placeholder gen; rename it to SPECIE_HUE
# and give every column a value of -1.
$ws.Range("A12").Value = "SPECIE_HUE"
$ws.Range("B12:L12").Value = -1

# Match the author's final selection in the saved file as closely as
# possible: the saved view has B12:L12 selected with C12 as the active
# cell (the user tabbed across the row while entering the -1 values).
$ws.Range("B12:L12").Select()
